$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns Q (billing_from) and R (billing_till).
# Format as Text first so the yyyy-mm-dd strings are stored literally
# instead of being auto-converted into date serial numbers.
$ws.Range("Q1:R57").NumberFormat = "@"

$ws.Range("Q1").Value = "billing_from"
$ws.Range("R1").Value = "billing_till"

$ws.Range("Q2").Value = '2025-02-01'
$ws.Range("R2").Value = '2025-02-27'
$ws.Range("Q3").Value = '2025-04-01'
$ws.Range("R3").Value = '2025-04-29'
$ws.Range("Q4").Value = '2025-05-01'
$ws.Range("R4").Value = '2025-05-30'
$ws.Range("Q5").Value = '2025-06-01'
$ws.Range("R5").Value = '2025-06-29'
$ws.Range("Q6").Value = '2025-08-01'
$ws.Range("R6").Value = '2025-08-30'
$ws.Range("Q7").Value = '2025-09-01'
$ws.Range("R7").Value = '2025-09-29'
$ws.Range("Q8").Value = '2025-11-01'
$ws.Range("R8").Value = '2025-11-29'
$ws.Range("Q9").Value = '2025-01-01'
$ws.Range("R9").Value = '2025-01-30'
$ws.Range("Q10").Value = '2025-02-01'
$ws.Range("R10").Value = '2025-02-27'
$ws.Range("Q11").Value = '2025-04-01'
$ws.Range("R11").Value = '2025-04-29'
$ws.Range("Q12").Value = '2025-05-01'
$ws.Range("R12").Value = '2025-05-30'
$ws.Range("Q13").Value = '2025-06-01'
$ws.Range("R13").Value = '2025-06-29'
$ws.Range("Q14").Value = '2025-08-01'
$ws.Range("R14").Value = '2025-08-30'
$ws.Range("Q15").Value = '2025-09-01'
$ws.Range("R15").Value = '2025-09-29'
$ws.Range("Q16").Value = '2025-11-01'
$ws.Range("R16").Value = '2025-11-29'
$ws.Range("Q17").Value = '2025-01-01'
$ws.Range("R17").Value = '2025-01-30'
$ws.Range("Q18").Value = '2025-02-01'
$ws.Range("R18").Value = '2025-02-27'
$ws.Range("Q19").Value = '2025-04-01'
$ws.Range("R19").Value = '2025-04-29'
$ws.Range("Q20").Value = '2025-05-01'
$ws.Range("R20").Value = '2025-05-30'
$ws.Range("Q21").Value = '2025-06-01'
$ws.Range("R21").Value = '2025-06-29'
$ws.Range("Q22").Value = '2025-01-01'
$ws.Range("R22").Value = '2025-01-31'
$ws.Range("Q23").Value = '2025-01-01'
$ws.Range("R23").Value = '2025-01-31'
$ws.Range("Q24").Value = '2025-01-01'
$ws.Range("R24").Value = '2025-01-31'
$ws.Range("Q25").Value = '2025-01-01'
$ws.Range("R25").Value = '2025-01-31'
$ws.Range("Q26").Value = '2025-01-01'
$ws.Range("R26").Value = '2025-01-31'
$ws.Range("Q27").Value = '2025-01-01'
$ws.Range("R27").Value = '2025-01-31'
$ws.Range("Q28").Value = '2025-01-01'
$ws.Range("R28").Value = '2025-01-31'
$ws.Range("Q29").Value = '2025-01-01'
$ws.Range("R29").Value = '2025-01-31'
$ws.Range("Q30").Value = '2025-01-01'
$ws.Range("R30").Value = '2025-01-31'
$ws.Range("Q31").Value = '2025-01-01'
$ws.Range("R31").Value = '2025-01-31'
$ws.Range("Q32").Value = '2025-01-01'
$ws.Range("R32").Value = '2025-01-31'
$ws.Range("Q33").Value = '2025-01-01'
$ws.Range("R33").Value = '2025-01-31'
$ws.Range("Q34").Value = '2025-06-01'
$ws.Range("R34").Value = '2025-06-30'
$ws.Range("Q35").Value = '2025-06-01'
$ws.Range("R35").Value = '2025-06-30'
$ws.Range("Q36").Value = '2025-06-01'
$ws.Range("R36").Value = '2025-06-30'
$ws.Range("Q37").Value = '2025-06-01'
$ws.Range("R37").Value = '2025-06-30'
$ws.Range("Q38").Value = '2025-06-01'
$ws.Range("R38").Value = '2025-06-30'
$ws.Range("Q39").Value = '2025-06-01'
$ws.Range("R39").Value = '2025-06-30'
$ws.Range("Q40").Value = '2025-06-01'
$ws.Range("R40").Value = '2025-06-30'
$ws.Range("Q41").Value = '2025-06-01'
$ws.Range("R41").Value = '2025-06-30'
$ws.Range("Q42").Value = '2025-06-01'
$ws.Range("R42").Value = '2025-06-30'
$ws.Range("Q43").Value = '2025-06-01'
$ws.Range("R43").Value = '2025-06-30'
$ws.Range("Q44").Value = '2025-06-01'
$ws.Range("R44").Value = '2025-06-30'
$ws.Range("Q45").Value = '2025-06-01'
$ws.Range("R45").Value = '2025-06-30'
$ws.Range("Q46").Value = '2025-06-01'
$ws.Range("R46").Value = '2025-06-30'
$ws.Range("Q47").Value = '2025-06-01'
$ws.Range("R47").Value = '2025-06-30'
$ws.Range("Q48").Value = '2025-06-01'
$ws.Range("R48").Value = '2025-06-30'
$ws.Range("Q49").Value = '2025-06-01'
$ws.Range("R49").Value = '2025-06-30'
$ws.Range("Q50").Value = '2025-06-01'
$ws.Range("R50").Value = '2025-06-30'
$ws.Range("Q51").Value = '2025-06-01'
$ws.Range("R51").Value = '2025-06-30'
$ws.Range("Q52").Value = '2025-06-01'
$ws.Range("R52").Value = '2025-06-30'
# Row 53: billing_from/billing_till left blank (no billing period data)
# Row 54: billing_from/billing_till left blank (no billing period data)
# Row 55: billing_from/billing_till left blank (no billing period data)
# Row 56: billing_from/billing_till left blank (no billing period data)
$ws.Range("Q57").Value = '2025-06-01'
$ws.Range("R57").Value = '2025-06-30'
